# Lecture 2.pptx — retitle the "space complexity" example slides.
#
# Slide 41 title: "Example I: The Factorial of a Number" -> "Space Complexity"
# Slide 42 title: "Example I: The Factorial of a Number" -> "Space Complexity: The Factorial Algorithm"

$p = $ppt.ActivePresentation

$slide41 = $p.Slides.Item(41)
$slide41.Shapes.Item(1).TextFrame.TextRange.Text = "Space Complexity"

$slide42 = $p.Slides.Item(42)
$slide42.Shapes.Item(1).TextFrame.TextRange.Text = "Space Complexity: The Factorial Algorithm"
